$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "dashboard and fullsearch" export adds two new product rows (23 and 24)
# at the bottom of the sheet (rows 13 and 14), extending the used range from
# A1:R12 to A1:R14. Columns left blank in the source export (B, E, G, I, M,
# N, O, P, Q, R for these rows) are simply left untouched/empty.

# --- Row 13 (product id 23 / "test") ---
$ws.Cells.Item(13, 1).Value = 23          # A13 id
$ws.Cells.Item(13, 3).Value = 239         # C13 barcode
$ws.Cells.Item(13, 4).Value = "test"      # D13 title
$ws.Cells.Item(13, 6).Value = 1           # F13 quantity
$ws.Cells.Item(13, 8).Value = 1.0         # H13 price

# J13/K13 (created_at/updated_at) need the same date number format as the
# rest of column J/K, so copy that formatting from an existing data row.
$ws.Cells.Item(13, 10).Value = 45283.87358929602
$ws.Cells.Item(2, 10).Copy()
$ws.Cells.Item(13, 10).PasteSpecial(-4122)
$ws.Cells.Item(13, 11).Value = 45283.88173452872
$ws.Cells.Item(2, 11).Copy()
$ws.Cells.Item(13, 11).PasteSpecial(-4122)

$ws.Cells.Item(13, 12).Value = 'http://localhost:3000/rails/active_storage/blobs/redirect/eyJfcmFpbHMiOnsibWVzc2FnZSI6IkJBaHBBcFVEIiwiZXhwIjpudWxsLCJwdXIiOiJibG9iX2lkIn19--fa7d7dcb19f78a8deb7ec8a7ba00f68ba73c2c0b/erp_fav.png http://localhost:3000/rails/active_storage/blobs/redirect/eyJfcmFpbHMiOnsibWVzc2FnZSI6IkJBaHBBcGNEIiwiZXhwIjpudWxsLCJwdXIiOiJibG9iX2lkIn19--034db0a4868dde326e1e9d8a47c66b5734f71fc3/Shopify%20Partners%202023-12-21%2019-38-52.png'  # L13 images

# --- Row 14 (product id 24 / "test2") ---
$ws.Cells.Item(14, 1).Value = 24          # A14 id
$ws.Cells.Item(14, 3).Value = 246         # C14 barcode
$ws.Cells.Item(14, 4).Value = "test2"     # D14 title
$ws.Cells.Item(14, 6).Value = 1           # F14 quantity
$ws.Cells.Item(14, 8).Value = 1.0         # H14 price

$ws.Cells.Item(14, 10).Value = 45283.88737998698
$ws.Cells.Item(2, 10).Copy()
$ws.Cells.Item(14, 10).PasteSpecial(-4122)
$ws.Cells.Item(14, 11).Value = 45283.888460052156
$ws.Cells.Item(2, 11).Copy()
$ws.Cells.Item(14, 11).PasteSpecial(-4122)

$ws.Cells.Item(14, 12).Value = 'http://localhost:3000/rails/active_storage/blobs/redirect/eyJfcmFpbHMiOnsibWVzc2FnZSI6IkJBaHBBcG9EIiwiZXhwIjpudWxsLCJwdXIiOiJibG9iX2lkIn19--732939eb1b2890b7bf2728c1621f14d7500d6e76/%D0%A2%D0%BE%D1%87%D0%BA%D0%B0%20%D0%A0%D0%BE%D1%81%D1%82%D0%B0%20-%20%D0%9F%D1%80%D0%BE%D0%B8%D0%B7%D0%B2%D0%BE%D0%BB%D1%8C%D0%BD%D1%8B%D0%B8%CC%86%20%D1%82%D0%BE%D0%B2%D0%B0%D1%80%20(vendor.model)%20-%20%D0%90%D0%B4%D0%B2%D0%B5%D0%BD%D1%82%D0%B5%D1%80%20-%20%D1%82%D0%BE%D0%B2%D0%B0%D1%80%D0%BD%D0%BE-%D0%BF%D1%80%D0%BE%D0%B8%D0%B7%D0%B2%D0%BE%D0%B4%D1%81%D1%82%D0%B2%D0%B5%D0%BD%D0%BD%D0%B0%D1%8F%20%D0%BA%D0%BE%D0%BC%D0%BF%D0%B0%D0%BD%D0%B8%D1%8F%20-%20InSales%202023-12-21%2012-14-53.png http://localhost:3000/rails/active_storage/blobs/redirect/eyJfcmFpbHMiOnsibWVzc2FnZSI6IkJBaHBBcHdEIiwiZXhwIjpudWxsLCJwdXIiOiJibG9iX2lkIn19--1682b57df3e50da59010c87fc126b5ea3654ad7e/Gmail%202023-12-13%2012-21-05.png'  # L14 images

$excel.CutCopyMode = $false

Write-Host "Added rows 13 and 14 (ids 23 and 24)"
